$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was added to the "Achicoria" table. It belongs
# chronologically before the existing row 113 (date 2021-06-09 / serial
# 44356), so insert a new row at 113 which pushes the existing rows 113-140
# down to 114-141, then populate the new row with its values.
$ws.Rows("113:113").Insert()

$ws.Range("A113").Value2 = 3
$ws.Range("B113").Value2 = "Femacal de La Calera"
$ws.Range("C113").Value2 = "Coquimbo"
$ws.Range("D113").Value2 = 44508
$ws.Range("E113").Value2 = 5
$ws.Range("F113").Value2 = 100112010
$ws.Range("G113").Value2 = "Achicoria"
$ws.Range("H113").Value2 = "Sin especificar"
$ws.Range("I113").Value2 = "Primera"
$ws.Range("J113").Value2 = 105
$ws.Range("K113").Value2 = 5800
$ws.Range("L113").Value2 = 6000
$ws.Range("M113").Value2 = 5895
$ws.Range("N113").Value2 = "$/caja 16 unidades"
$ws.Range("O113").Value2 = "Provincia de Quillota"
$ws.Range("P113").Value2 = 368
$ws.Range("Q113").Value2 = 16
$ws.Range("R113").Value2 = "Hortaliza"
